$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct a handful of previously-recorded QAM2 values ---
$ws.Range("C146:F146").Value = 562328000000
$ws.Range("C149:F149").Value = 557944000000
$ws.Range("C150:F150").Value = 554911000000
$ws.Range("C152:F152").Value = 554033000000
$ws.Range("C153:F153").Value = 551532000000
$ws.Range("C200:F200").Value = 683766000000

# --- Append the two newest monthly QAM2 data points ---
# Clone row 200's formatting (date style, borders, alignment) onto the new rows first
$ws.Range("A200:G200").Copy($ws.Range("A201:G201"))
$ws.Range("A200:G200").Copy($ws.Range("A202:G202"))

$ws.Range("A201").Value = 45139.41666666666
$ws.Range("B201").Value = "ECONOMICS:QAM2"
$ws.Range("C201:F201").Value = 685009000000
$ws.Range("G201").Value = 0

$ws.Range("A202").Value = 45170.41666666666
$ws.Range("B202").Value = "ECONOMICS:QAM2"
$ws.Range("C202:F202").Value = 702188000000
$ws.Range("G202").Value = 0
